$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F84").Value = 26
$ws.Range("G84").Value = 2475.46
$ws.Range("F86").Value = 43
$ws.Range("G86").Value = 4852.12
$ws.Range("B114").Value = 302893.86
$ws.Range("F150").Value = 294
$ws.Range("G150").Value = 28462.14
$ws.Range("B152").Value = 23815.7
$ws.Range("F244").Value = 30
$ws.Range("G244").Value = 999
$ws.Range("B246").Value = 48706
$ws.Range("E246").Value = 39.8
$ws.Range("F246").Value = -144
$ws.Range("G246").Value = -4795.2
$ws.Range("B247").Value = 64973
$ws.Range("E247").Value = 35.4
$ws.Range("F247").Value = 114
$ws.Range("G247").Value = 3796.2
$ws.Range("B274").Value = 113542.67
$ws.Range("B277").Value = 61610
$ws.Range("E277").Value = 122.71
$ws.Range("F277").Value = -58
$ws.Range("G277").Value = -5957.18
$ws.Range("B278").Value = 63565
$ws.Range("E278").Value = 109.19
$ws.Range("F278").Value = 60
$ws.Range("G278").Value = 6162.6
$ws.Range("B294").Value = 57802
$ws.Range("E294").Value = 162.71
$ws.Range("F294").Value = -79
$ws.Range("G294").Value = -11334.92
$ws.Range("B295").Value = 63571
$ws.Range("F295").Value = 9
$ws.Range("G295").Value = 1291.32
$ws.Range("B296").Value = 63531
$ws.Range("E296").Value = 152.53
$ws.Range("F296").Value = 80
$ws.Range("G296").Value = 11478.4
$ws.Range("B299").Value = 63510
$ws.Range("E299").Value = 50.66
$ws.Range("F299").Value = 150
$ws.Range("G299").Value = 7146
$ws.Range("B300").Value = 55356
$ws.Range("E300").Value = 54.04
$ws.Range("F300").Value = -158
$ws.Range("G300").Value = -7527.12
$ws.Range("B311").Value = 63563
$ws.Range("E311").Value = 119.04
$ws.Range("F311").Value = 2
$ws.Range("G311").Value = 223.92
$ws.Range("B312").Value = 61605
$ws.Range("E312").Value = 133.78
$ws.Range("F312").Value = -13
$ws.Range("G312").Value = -1455.48
$ws.Range("F314").Value = 6
$ws.Range("G314").Value = 774.0599999999999
$ws.Range("F321").Value = 164
$ws.Range("G321").Value = 9626.799999999999
$ws.Range("F328").Value = 1406
$ws.Range("G328").Value = 29568.18
$ws.Range("F329").Value = 53
$ws.Range("G329").Value = 8533
$ws.Range("B339").Value = 383110.55
$ws.Range("F355").Value = 121
$ws.Range("G355").Value = 3901.04
$ws.Range("B356").Value = 63681
$ws.Range("E356").Value = 23.84
$ws.Range("F356").Value = 0
$ws.Range("G356").Value = 0
$ws.Range("B357").Value = 31930
$ws.Range("E357").Value = 26.8
$ws.Range("F357").Value = -62
$ws.Range("G357").Value = -1390.04
$ws.Range("B361").Value = 14067.36
$ws.Range("F389").Value = 1
$ws.Range("G389").Value = 55.86
$ws.Range("B395").Value = 271464.89
$ws.Range("B420").Value = 47097
$ws.Range("D420").Value = 112.28
$ws.Range("E420").Value = 134.16
$ws.Range("F420").Value = 15
$ws.Range("G420").Value = 1684.2
$ws.Range("B421").Value = 58047
$ws.Range("D421").Value = 105.54
$ws.Range("E421").Value = 126.1
$ws.Range("F421").Value = 43
$ws.Range("G421").Value = 4538.22
$ws.Range("B472").Value = 64915
$ws.Range("E472").Value = 20.98
$ws.Range("F472").Value = 0
$ws.Range("G472").Value = 0
$ws.Range("B473").Value = 45695
$ws.Range("E473").Value = 23.58
$ws.Range("F473").Value = -36
$ws.Range("G473").Value = -710.28
$ws.Range("F477").Value = 144
$ws.Range("G477").Value = 2841.12
$ws.Range("B479").Value = 45718
$ws.Range("E479").Value = 19.38
$ws.Range("F479").Value = -294
$ws.Range("G479").Value = -4768.68
$ws.Range("B480").Value = 64927
$ws.Range("E480").Value = 17.26
$ws.Range("F480").Value = 227
$ws.Range("G480").Value = 3681.94
$ws.Range("F484").Value = 679
$ws.Range("G484").Value = 4406.71
$ws.Range("F485").Value = 236
$ws.Range("G485").Value = 3103.4
$ws.Range("B492").Value = 4242.3
$ws.Range("F497").Value = 52
$ws.Range("G497").Value = 2604.16
$ws.Range("B508").Value = 20970.48
$ws.Range("F511").Value = 26
$ws.Range("G511").Value = 630.5
$ws.Range("B528").Value = 21004.38
$ws.Range("F648").Value = 91
$ws.Range("G648").Value = 9458.540000000001
$ws.Range("B651").Value = 40412.53
$ws.Range("F701").Value = 234
$ws.Range("G701").Value = 33492.42
$ws.Range("F705").Value = 96
$ws.Range("G705").Value = 7265.28
$ws.Range("F707").Value = 169
$ws.Range("G707").Value = 3670.68
$ws.Range("F708").Value = 69
$ws.Range("G708").Value = 2572.32
$ws.Range("F714").Value = 42
$ws.Range("G714").Value = 1571.64
$ws.Range("B716").Value = 211268.81
$ws.Range("F720").Value = 73
$ws.Range("G720").Value = 11963.97
$ws.Range("F724").Value = 57
$ws.Range("G724").Value = 8573.370000000001
$ws.Range("F725").Value = 17
$ws.Range("G725").Value = 2487.61
$ws.Range("F733").Value = 71
$ws.Range("G733").Value = 2901.77
$ws.Range("F742").Value = 135
$ws.Range("G742").Value = 7693.65
$ws.Range("B743").Value = 115439.4
$ws.Range("F771").Value = 529
$ws.Range("G771").Value = 76519.85000000001
$ws.Range("B775").Value = 890418.88
$ws.Range("B793").Value = 3528564.2
$ws.Range("B794").Value = 3528564.2
